$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet: excess_return_without_cost
$ws1.Range("C2").Value = 0.0012
$ws1.Range("C3").Value = 0.00228
$ws1.Range("D3").Value = 0.00231
$ws1.Range("C4").Value = 0.00115
$ws1.Range("D4").Value = 0.00135
$ws1.Range("E4").Value = 0.00167
$ws1.Range("F4").Value = 0.00174
$ws1.Range("C5").Value = 0.00079
$ws1.Range("D5").Value = 0.0009700000000000001
$ws1.Range("E5").Value = 0.00153
$ws1.Range("F5").Value = 0.00164
$ws1.Range("G5").Value = 0.0016
$ws1.Range("C6").Value = 0.00056
$ws1.Range("D6").Value = 0.00079
$ws1.Range("E6").Value = 0.00125
$ws1.Range("F6").Value = 0.0017
$ws1.Range("G6").Value = 0.00174
$ws1.Range("C7").Value = 0.00056
$ws1.Range("D7").Value = 0.00074
$ws1.Range("E7").Value = 0.00138
$ws1.Range("F7").Value = 0.00152
$ws1.Range("G7").Value = 0.00162
$ws1.Range("C8").Value = 0.01794
$ws1.Range("C9").Value = 0.01906
$ws1.Range("D9").Value = 0.01775
$ws1.Range("C10").Value = 0.01316
$ws1.Range("D10").Value = 0.01382
$ws1.Range("E10").Value = 0.01344
$ws1.Range("F10").Value = 0.01329
$ws1.Range("C11").Value = 0.01118
$ws1.Range("D11").Value = 0.0113
$ws1.Range("E11").Value = 0.01163
$ws1.Range("F11").Value = 0.01164
$ws1.Range("G11").Value = 0.01151
$ws1.Range("C12").Value = 0.00979
$ws1.Range("D12").Value = 0.01035
$ws1.Range("E12").Value = 0.01029
$ws1.Range("F12").Value = 0.01054
$ws1.Range("G12").Value = 0.0104
$ws1.Range("C13").Value = 0.009480000000000001
$ws1.Range("D13").Value = 0.00957
$ws1.Range("E13").Value = 0.009730000000000001
$ws1.Range("F13").Value = 0.009549999999999999
$ws1.Range("G13").Value = 0.009549999999999999
$ws1.Range("C14").Value = 0.28451
$ws1.Range("C15").Value = 0.54167
$ws1.Range("D15").Value = 0.5502899999999999
$ws1.Range("C16").Value = 0.27434
$ws1.Range("D16").Value = 0.32113
$ws1.Range("E16").Value = 0.3982
$ws1.Range("F16").Value = 0.41333
$ws1.Range("C17").Value = 0.18696
$ws1.Range("D17").Value = 0.23112
$ws1.Range("E17").Value = 0.36317
$ws1.Range("F17").Value = 0.39128
$ws1.Range("G17").Value = 0.38152
$ws1.Range("C18").Value = 0.13296
$ws1.Range("D18").Value = 0.18822
$ws1.Range("E18").Value = 0.29635
$ws1.Range("F18").Value = 0.4053
$ws1.Range("G18").Value = 0.41498
$ws1.Range("C19").Value = 0.13249
$ws1.Range("D19").Value = 0.17518
$ws1.Range("E19").Value = 0.32846
$ws1.Range("F19").Value = 0.36066
$ws1.Range("G19").Value = 0.38533
$ws1.Range("C20").Value = 1.02791
$ws1.Range("C21").Value = 1.84239
$ws1.Range("D21").Value = 2.01014
$ws1.Range("C22").Value = 1.35169
$ws1.Range("D22").Value = 1.50659
$ws1.Range("E22").Value = 1.92108
$ws1.Range("F22").Value = 2.01618
$ws1.Range("C23").Value = 1.08418
$ws1.Range("D23").Value = 1.32552
$ws1.Range("E23").Value = 2.0249
$ws1.Range("F23").Value = 2.17832
$ws1.Range("G23").Value = 2.14786
$ws1.Range("C24").Value = 0.8801099999999999
$ws1.Range("D24").Value = 1.17933
$ws1.Range("E24").Value = 1.86632
$ws1.Range("F24").Value = 2.49304
$ws1.Range("G24").Value = 2.5871
$ws1.Range("C25").Value = 0.90626
$ws1.Range("D25").Value = 1.186
$ws1.Range("E25").Value = 2.18774
$ws1.Range("F25").Value = 2.44809
$ws1.Range("G25").Value = 2.61492
$ws1.Range("C26").Value = -0.37933
$ws1.Range("C27").Value = -0.27995
$ws1.Range("D27").Value = -0.25559
$ws1.Range("C28").Value = -0.23494
$ws1.Range("D28").Value = -0.23298
$ws1.Range("E28").Value = -0.16415
$ws1.Range("F28").Value = -0.13614
$ws1.Range("C29").Value = -0.21008
$ws1.Range("D29").Value = -0.17116
$ws1.Range("E29").Value = -0.11868
$ws1.Range("F29").Value = -0.11145
$ws1.Range("G29").Value = -0.11952
$ws1.Range("C30").Value = -0.20593
$ws1.Range("D30").Value = -0.17864
$ws1.Range("E30").Value = -0.10728
$ws1.Range("F30").Value = -0.11495
$ws1.Range("G30").Value = -0.09958
$ws1.Range("C31").Value = -0.24809
$ws1.Range("D31").Value = -0.18943
$ws1.Range("E31").Value = -0.10674
$ws1.Range("F31").Value = -0.09878000000000001
$ws1.Range("G31").Value = -0.10198

# Sheet: excess_return_with_cost
$ws2.Range("C2").Value = 0.00061
$ws2.Range("C3").Value = 0.00128
$ws2.Range("D3").Value = 0.00114
$ws2.Range("C4").Value = 0.00066
$ws2.Range("D4").Value = 0.00051
$ws2.Range("E4").Value = 0.00065
$ws2.Range("F4").Value = 0.00065
$ws2.Range("C5").Value = 0.00043
$ws2.Range("D5").Value = 0.00035
$ws2.Range("E5").Value = 0.00066
$ws2.Range("F5").Value = 0.00066
$ws2.Range("G5").Value = 0.00058
$ws2.Range("C6").Value = 0.00029
$ws2.Range("D6").Value = 0.00028
$ws2.Range("E6").Value = 0.00054
$ws2.Range("F6").Value = 0.0008
$ws2.Range("G6").Value = 0.00077
$ws2.Range("C7").Value = 0.00033
$ws2.Range("D7").Value = 0.00033
$ws2.Range("E7").Value = 0.00077
$ws2.Range("F7").Value = 0.00074
$ws2.Range("G7").Value = 0.00075
$ws2.Range("C8").Value = 0.01792
$ws2.Range("C9").Value = 0.01907
$ws2.Range("D9").Value = 0.01775
$ws2.Range("C10").Value = 0.01315
$ws2.Range("D10").Value = 0.01381
$ws2.Range("E10").Value = 0.01342
$ws2.Range("F10").Value = 0.01327
$ws2.Range("C11").Value = 0.01118
$ws2.Range("D11").Value = 0.0113
$ws2.Range("E11").Value = 0.01162
$ws2.Range("F11").Value = 0.01164
$ws2.Range("G11").Value = 0.0115
$ws2.Range("C12").Value = 0.0098
$ws2.Range("D12").Value = 0.01035
$ws2.Range("E12").Value = 0.0103
$ws2.Range("F12").Value = 0.01054
$ws2.Range("G12").Value = 0.0104
$ws2.Range("C13").Value = 0.009480000000000001
$ws2.Range("D13").Value = 0.00957
$ws2.Range("E13").Value = 0.009730000000000001
$ws2.Range("F13").Value = 0.009549999999999999
$ws2.Range("G13").Value = 0.009560000000000001
$ws2.Range("C14").Value = 0.14545
$ws2.Range("C15").Value = 0.30422
$ws2.Range("D15").Value = 0.27104
$ws2.Range("C16").Value = 0.15686
$ws2.Range("D16").Value = 0.1208
$ws2.Range("E16").Value = 0.15357
$ws2.Range("F16").Value = 0.15528
$ws2.Range("C17").Value = 0.10251
$ws2.Range("D17").Value = 0.083
$ws2.Range("E17").Value = 0.15796
$ws2.Range("F17").Value = 0.15669
$ws2.Range("G17").Value = 0.13889
$ws2.Range("C18").Value = 0.06970999999999999
$ws2.Range("D18").Value = 0.06753000000000001
$ws2.Range("E18").Value = 0.12825
$ws2.Range("F18").Value = 0.1901
$ws2.Range("G18").Value = 0.1841
$ws2.Range("C19").Value = 0.07940999999999999
$ws2.Range("D19").Value = 0.07947
$ws2.Range("E19").Value = 0.18235
$ws2.Range("F19").Value = 0.17682
$ws2.Range("G19").Value = 0.17769
$ws2.Range("C20").Value = 0.52604
$ws2.Range("C21").Value = 1.03422
$ws2.Range("D21").Value = 0.98961
$ws2.Range("C22").Value = 0.7730399999999999
$ws2.Range("D22").Value = 0.56714
$ws2.Range("E22").Value = 0.742
$ws2.Range("F22").Value = 0.75874
$ws2.Range("C23").Value = 0.59433
$ws2.Range("D23").Value = 0.47615
$ws2.Range("E23").Value = 0.88088
$ws2.Range("F23").Value = 0.8728900000000001
$ws2.Range("G23").Value = 0.78254
$ws2.Range("C24").Value = 0.46128
$ws2.Range("D24").Value = 0.42311
$ws2.Range("E24").Value = 0.80736
$ws2.Range("F24").Value = 1.16873
$ws2.Range("G24").Value = 1.14782
$ws2.Range("C25").Value = 0.54313
$ws2.Range("D25").Value = 0.53806
$ws2.Range("E25").Value = 1.21448
$ws2.Range("F25").Value = 1.19959
$ws2.Range("G25").Value = 1.20509
$ws2.Range("C26").Value = -0.4446
$ws2.Range("C27").Value = -0.36801
$ws2.Range("D27").Value = -0.33815
$ws2.Range("C28").Value = -0.32277
$ws2.Range("D28").Value = -0.33395
$ws2.Range("E28").Value = -0.26302
$ws2.Range("F28").Value = -0.24434
$ws2.Range("C29").Value = -0.32437
$ws2.Range("D29").Value = -0.27803
$ws2.Range("E29").Value = -0.19298
$ws2.Range("F29").Value = -0.20474
$ws2.Range("G29").Value = -0.20455
$ws2.Range("C30").Value = -0.27065
$ws2.Range("D30").Value = -0.29984
$ws2.Range("E30").Value = -0.16614
$ws2.Range("F30").Value = -0.15213
$ws2.Range("G30").Value = -0.17776
$ws2.Range("C31").Value = -0.30084
$ws2.Range("D31").Value = -0.29343
$ws2.Range("E31").Value = -0.13976
$ws2.Range("F31").Value = -0.12559
$ws2.Range("G31").Value = -0.12936

Write-Output "updated report values"